$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.00121430506933784
$ws.Range("D2").Value = 0.00204646281340537
$ws.Range("E2").Value = 0.4304973541417638
$ws.Range("F2").Value = 0.8133486523056916
$ws.Range("G2").Value = 0.002341840673280299
$ws.Range("I2").Value = 0.6112125958755712
$ws.Range("N2").Value = 2.873581687696344
$ws.Range("O2").Value = 2.71929717795183
$ws.Range("C3").Value = 0.001082874213395257
$ws.Range("D3").Value = 0.001855303547857901
$ws.Range("E3").Value = 0.3750918814098725
$ws.Range("F3").Value = 0.7525647962886808
$ws.Range("G3").Value = 0.002346712352590319
$ws.Range("I3").Value = 0.565059335681795
$ws.Range("N3").Value = 2.562605684679454
$ws.Range("O3").Value = 2.515198469554434
$ws.Range("C4").Value = 0.001002524826834161
$ws.Range("D4").Value = 0.001739784820802726
$ws.Range("E4").Value = 0.3412096363495465
$ws.Range("F4").Value = 0.7157517958633264
$ws.Range("G4").Value = 0.002349859887713835
$ws.Range("I4").Value = 0.5371031260879988
$ws.Range("N4").Value = 2.371325805375761
$ws.Range("O4").Value = 2.391593294181121
$ws.Range("C5").Value = 0.0009698708156626878
$ws.Range("D5").Value = 0.001693152239583995
$ws.Range("E5").Value = 0.3274337357478458
$ws.Range("F5").Value = 0.7008764035240631
$ws.Range("G5").Value = 0.002351181983377026
$ws.Range("I5").Value = 0.5258055262832073
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 2.341647831043247
$ws.Range("C6").Value = 0.0009644540519069267
$ws.Range("D6").Value = 0.001685434953795806
$ws.Range("E6").Value = 0.3251480654236474
$ws.Range("F6").Value = 0.6984139299659091
$ws.Range("G6").Value = 0.002351403903144679
$ws.Range("I6").Value = 0.5239352555636714
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 2.333379901361866
$ws.Range("C7").Value = 0.001002084080932164
$ws.Range("D7").Value = 0.001739154157634815
$ws.Range("E7").Value = 0.3410237266118799
$ws.Range("F7").Value = 0.7155506723287317
$ws.Range("G7").Value = 0.002349877557992056
$ws.Range("I7").Value = 0.5369503806425655
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 2.390918000611521
$ws.Range("C8").Value = 0.001168915712886687
$ws.Range("D8").Value = 0.001980151755827109
$ws.Range("E8").Value = 0.4113632681818302
$ws.Range("F8").Value = 0.7922838023180532
$ws.Range("G8").Value = 0.002343488075216054
$ws.Range("I8").Value = 0.5952188533536429
$ws.Range("N8").Value = 2.766433886209654
$ws.Range("O8").Value = 2.648564921756588
$ws.Range("C9").Value = 0.001498820731374195
$ws.Range("D9").Value = 0.00246857337649331
$ws.Range("E9").Value = 0.5505299427434807
$ws.Range("F9").Value = 0.9468763800419708
$ws.Range("G9").Value = 0.002332191884470869
$ws.Range("I9").Value = 0.7125798299844917
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 3.16769198211972
$ws.Range("C10").Value = 0.001742874978909015
$ws.Range("D10").Value = 0.002838666599469519
$ws.Range("E10").Value = 0.6537366696595939
$ws.Range("F10").Value = 1.06309597156465
$ws.Range("G10").Value = 0.002324635282823703
$ws.Range("I10").Value = 0.8007924982839825
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 3.558010637400116
$ws.Range("C11").Value = 0.001854267906850282
$ws.Range("D11").Value = 0.003009832214619479
$ws.Range("E11").Value = 0.7009434488239918
$ws.Range("F11").Value = 1.116569624944674
$ws.Range("G11").Value = 0.00232135686987699
$ws.Range("I11").Value = 0.8413766745340467
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 3.737614285049801
$ws.Range("C12").Value = 0.001896502724171967
$ws.Range("D12").Value = 0.003075081495182985
$ws.Range("E12").Value = 0.7188601408737298
$ws.Range("F12").Value = 1.136907811103981
$ws.Range("G12").Value = 0.002320138148584871
$ws.Range("I12").Value = 0.856812041308217
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 3.805927218431975
$ws.Range("C13").Value = 0.001887404365003675
$ws.Range("D13").Value = 0.003061009191064556
$ws.Range("E13").Value = 0.7149996029872341
$ws.Range("F13").Value = 1.132523635391408
$ws.Range("G13").Value = 0.002320399612628563
$ws.Range("I13").Value = 0.8534847543662067
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 3.791201313206784
$ws.Range("C14").Value = 0.00185774153795748
$ws.Range("D14").Value = 0.003015191475455481
$ws.Range("E14").Value = 0.7024166338525077
$ws.Range("F14").Value = 1.118241067548041
$ws.Range("G14").Value = 0.00232125614985158
$ws.Range("I14").Value = 0.8426451997138429
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 3.743228361005436
$ws.Range("C15").Value = 0.001839579047743456
$ws.Range("D15").Value = 0.002987184001156606
$ws.Range("E15").Value = 0.694714581711608
$ws.Range("F15").Value = 1.10950421354994
$ws.Range("G15").Value = 0.002321783761383234
$ws.Range("I15").Value = 0.8360144317726395
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 3.713882935004619
$ws.Range("C16").Value = 0.001735602602920494
$ws.Range("D16").Value = 0.002827539476296437
$ws.Range("E16").Value = 0.6506570735473218
$ws.Range("F16").Value = 1.059613674597244
$ws.Range("G16").Value = 0.002324852725070079
$ws.Range("I16").Value = 0.798149524367247
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 3.546314856958759
$ws.Range("C17").Value = 0.00167191103755826
$ws.Range("D17").Value = 0.002730342243449257
$ws.Range("E17").Value = 0.623697632425575
$ws.Range("F17").Value = 1.029163690499303
$ws.Range("G17").Value = 0.00232677609088879
$ws.Range("I17").Value = 0.7750384128584074
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 3.444046001339359
$ws.Range("C18").Value = 0.001635312328900795
$ws.Range("D18").Value = 0.002674699540431646
$ws.Range("E18").Value = 0.6082153024572534
$ws.Range("F18").Value = 1.011706473618517
$ws.Range("G18").Value = 0.002327897345805218
$ws.Range("I18").Value = 0.7617883241393457
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 3.385415803845092
$ws.Range("C19").Value = 0.001622926667021574
$ws.Range("D19").Value = 0.002655904109918339
$ws.Range("E19").Value = 0.6029772733294863
$ws.Range("F19").Value = 1.005805459256663
$ws.Range("G19").Value = 0.002328279561517547
$ws.Range("I19").Value = 0.7573093800281043
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 3.365597423163877
$ws.Range("C20").Value = 0.001678687496848852
$ws.Range("D20").Value = 0.002740661673289679
$ws.Range("E20").Value = 0.6265649991965745
$ws.Range("F20").Value = 1.032399250097257
$ws.Range("G20").Value = 0.002326569795315775
$ws.Range("I20").Value = 0.7774941886937512
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("O20").Value = 3.454912760925026
$ws.Range("C21").Value = 0.001866452803493956
$ws.Range("D21").Value = 0.003028637277004975
$ws.Range("E21").Value = 0.7061114294735802
$ws.Range("F21").Value = 1.122433776554914
$ws.Range("G21").Value = 0.0023210039480104
$ws.Range("I21").Value = 0.8458272088643355
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("O21").Value = 3.757310958386086
$ws.Range("C22").Value = 0.001989475413587627
$ws.Range("D22").Value = 0.003219381215664896
$ws.Range("E22").Value = 0.7583373973780283
$ws.Range("F22").Value = 1.181795479378252
$ws.Range("G22").Value = 0.002317498841747085
$ws.Range("I22").Value = 0.8908781354258508
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("O22").Value = 3.956702928797881
$ws.Range("C23").Value = 0.001923788026378048
$ws.Range("D23").Value = 0.003117335930564735
$ws.Range("E23").Value = 0.7304404936382269
$ws.Range("F23").Value = 1.150064900034465
$ws.Range("G23").Value = 0.002319357505304813
$ws.Range("I23").Value = 0.8667973049203681
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 3.850120629423031
$ws.Range("C24").Value = 0.001675623802238846
$ws.Range("D24").Value = 0.002735995522815671
$ws.Range("E24").Value = 0.6252686099535509
$ws.Range("F24").Value = 1.030936301482228
$ws.Range("G24").Value = 0.002326663013365198
$ws.Range("I24").Value = 0.7763838179350273
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 3.449999383836655
$ws.Range("C25").Value = 0.001409279053220303
$ws.Range("D25").Value = 0.002334584835804066
$ws.Range("E25").Value = 0.5127283134970213
$ws.Range("F25").Value = 0.9046007503701361
$ws.Range("G25").Value = 0.002335116704995537
$ws.Range("I25").Value = 0.6804888268093805
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 3.025721714252711
